$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5484259724617004
$ws.Range("B1").Value = 1.212943434715271
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.087329864501953
$ws.Range("E1").Value = 1.108350396156311
